# Weekly Fruta/Hortaliza update:
# Insert 3 new rows (new week of data, 2021-11-22) above the existing
# row 105, pushing the old rows 105-119 down to 108-122.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 105 (shifts existing rows 105-119 down to 108-122)
$ws.Range("A105:A107").EntireRow.Insert()

# New data for the week of 2021-11-22 ("Provincia de Limarí")
$newRows = @(
    @{ Row = 105; D = [DateTime]"2021-11-22"; L = "Especial"; M = 300; N = 1800; O = 1900; P = 1850; S = 1850; T = 1 },
    @{ Row = 106; D = [DateTime]"2021-11-22"; L = "Primera";  M = 300; N = 1500; O = 1600; P = 1550; S = 1550; T = 1 },
    @{ Row = 107; D = [DateTime]"2021-11-22"; L = "Segunda";  M = 240; N = 1200; O = 1300; P = 1250; S = 1250; T = 1 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $r.D
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/kilo (en caja de 15 kilos)"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
